{"js": "// MLD: drop \"DateNaissance\" from the USER row's attribute list\n// (\", Pseudo, Email, Password, DateNaissance, Salt)\" ->\n//  \", Pseudo, Email, Password, Salt)\").\n// Word also re-tracks its \"_GoBack\" (last-edit) bookmark to the spot of\n// this edit, so move that bookmark from the trailing empty paragraph to\n// right after \"Password\".\n\nconst body = context.document.body;\n\n// 1) Drop the old _GoBack bookmark (lived on the last, empty paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Remove \", DateNaissance\" from the USER row definition.\nconst dateNaissance = body.search(\", DateNaissance\", { matchCase: true });\ndateNaissance.load(\"text\");\nawait context.sync();\n\nif (dateNaissance.items.length > 0) {\n  dateNaissance.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Re-create _GoBack right after \"Password\" - where the text was edited.\nconst password = body.search(\"Password\", { matchCase: true });\npassword.load(\"text\");\nawait context.sync();\n\nif (password.items.length > 0) {\n  const afterPassword = password.items[0].getRange(Word.RangeLocation.after);\n  afterPassword.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# MLD: drop \"DateNaissance\" from the USER row's attribute list\n# (\", Pseudo, Email, Password, DateNaissance, Salt)\" ->\n#  \", Pseudo, Email, Password, Salt)\").\n# Word also re-tracks its \"_GoBack\" (last-edit) bookmark to the spot of\n# this edit, so move that bookmark from the trailing empty paragraph to\n# right after \"Password\".\n\n$d = $word.ActiveDocument\n\n# 1) Remove \", DateNaissance\" from the USER row definition.\n$find = $d.Content.Find\n$find.Text = \", DateNaissance\"\n$find.Replacement.Text = \"\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) Re-create _GoBack right after \"Password\" - where the text was edited.\n#    Bookmarks.Add with the existing \"_GoBack\" name moves it (removing the\n#    old one on the trailing empty paragraph).\n$rng = $d.Content\n$find2 = $rng.Find\n$find2.Text = \"Password\"\n$find2.Execute() | Out-Null\n$collapsed = $d.Range($rng.End, $rng.End)\n$d.Bookmarks.Add(\"_GoBack\", $collapsed) | Out-Null\n"}
